# Applies the commit's changes to the workbook:
#  - Stats sheet: becomes the active/selected sheet again (tabSelected),
#    selection moves to E29, and the formulas that pull weekly totals from
#    Journal are bumped down by two rows (because Journal gains two rows).
#  - Journal sheet: the "Semaine 3" block gets three new work-log entries
#    (previously just an empty placeholder), which pushes every section
#    below it down by two rows; two extra filler rows are appended at the
#    bottom to keep the sheet's visual length, and the sheet view / active
#    tab move from Journal to Stats.
#  - Workbook: no sheet is forced "active" on open any more.

$wb = $excel.ActiveWorkbook

$stats = $wb.Worksheets.Item("Stats")
$journal = $wb.Worksheets.Item("Journal")

# --- Journal: insert the 2 extra rows the new entries need -----------------
# Row 44 is already a blank "data" row for the Semaine 3 block and row 45 is
# the blank separator row above the "Total /sem" line; inserting two rows at
# 45 keeps row 44 in place and shifts everything from the old row 45 onward
# down by two (Excel auto-adjusts the SUM()/formula ranges as it goes).
$journal.Rows.Item(45).Insert()
$journal.Rows.Item(45).Insert()

# --- Journal: give the 3 new data rows (44:46) the same look as the other
# data rows in this block (row 36 is a normal, non-customFormat data row) ---
$journal.Range("A36:E36").Copy()
$journal.Range("A44:E46").PasteSpecial(-4122)
$journal.Rows.Item(44).OutlineLevel = 1
$journal.Rows.Item(45).OutlineLevel = 1
$journal.Rows.Item(46).OutlineLevel = 1

# --- Journal: fill in the 3 new work-log entries ----------------------------
$journal.Cells.Item(44, 1).Value = 43600
$journal.Cells.Item(44, 2).Value = 3
$journal.Cells.Item(44, 3).Value = 0.041666666666666664
$journal.Cells.Item(44, 4).Value = "Réalisation"
$journal.Cells.Item(44, 5).Value = "Mise en place des options - intégration du type de giratoire"

$journal.Cells.Item(45, 1).Value = 43600
$journal.Cells.Item(45, 2).Value = 3
$journal.Cells.Item(45, 3).Value = 0.08333333333333333
$journal.Cells.Item(45, 4).Value = "Réalisation"
$journal.Cells.Item(45, 5).Value = "Mise en place des options - intégration du nombre de routes"

$journal.Cells.Item(46, 1).Value = 43600
$journal.Cells.Item(46, 2).Value = 3
$journal.Cells.Item(46, 3).Value = 0.03125
$journal.Cells.Item(46, 4).Value = "Réalisation"
$journal.Cells.Item(46, 5).Value = "Mise en place des options - intégration du nombre de véhicules par route"

# --- Journal: two extra filler rows at the end (65:66), matching 61:64 -----
# (column B is skipped on purpose - row 64 has no B cell either)
$journal.Range("A64").Copy()
$journal.Range("A65:A66").PasteSpecial(-4122)
$journal.Range("C64:F64").Copy()
$journal.Range("C65:F66").PasteSpecial(-4122)

# --- Journal: sheet view moves off the Journal tab, scrolled & selected ----
$journal.Range("C47").Select()

# --- Stats: formulas shift down two rows (C46->C48, C50->C52, ...) --------
$stats.Range("B2").Formula = "=Journal!C48"
$stats.Range("B3").Formula = "=Journal!C52"
$stats.Range("B4").Formula = "=Journal!C56"
$stats.Range("B5").Formula = "=Journal!C60"

# --- Stats becomes the selected/active sheet, with E29 selected -----------
$stats.Activate()
$stats.Range("E29").Select()

$wb.Save()
